# Commit message: Added M0 and M1 to tasks sheet. Populated/updated M0
#
# Adds two new worksheets ('M0 - Account Mgmt' and 'M1 - Game Data') between
# 'M-0.5' and 'Links', populates M0 with the account-mgmt task breakdown,
# lays out the (currently empty) section headers on M1, and makes M0 the
# active/selected sheet, matching the target workbook.

$wb = $excel.ActiveWorkbook

# ---- Create the two new sheets in position, right after 'M-0.5' ----
$afterM05 = $wb.Worksheets.Item('M-0.5')
$wsM0 = $wb.Worksheets.Add($null, $afterM05)
$wsM0.Name = 'M0 - Account Mgmt'

$afterM0 = $wb.Worksheets.Item('M0 - Account Mgmt')
$wsM1 = $wb.Worksheets.Add($null, $afterM0)
$wsM1.Name = 'M1 - Game Data'

# Reference sheet that already carries the 3 coloured section-header styles
# (s=1 gold/teal 'MISC', s=2 blue 'ADMIN/APPLICATION', s=3 green 'USER/MISC TASKS')
# that we need to reapply on the new sheets - reuse them via copy/paste-format
# instead of guessing raw style indices.
$styleSrc = $wb.Worksheets.Item('M-1 Tasks')

# ===================== M0 - Account Mgmt =====================
$ws = $wsM0

# Column widths
$ws.Columns.Item(1).ColumnWidth = 54.0
$ws.Columns.Item(2).ColumnWidth = 14.333333333333334
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666
$ws.Columns.Item(4).ColumnWidth = 71.83333333333333

# Cell text
$ws.Range('A1').Value = 'ADMIN FEATURES'
$ws.Range('B1').Value = 'Area / Component'
$ws.Range('C1').Value = 'Status'
$ws.Range('D1').Value = 'Notes'
$ws.Range('A2').Value = 'Create Admin View (and model/controller if required)'
$ws.Range('B2').Value = 'Admin Controller'
$ws.Range('A3').Value = 'Change Admin password'
$ws.Range('B3').Value = 'Admin Controller'
$ws.Range('C3').Value = 'Done'
$ws.Range('A4').Value = 'Change User password'
$ws.Range('B4').Value = 'Admin Controller'
$ws.Range('A5').Value = 'Delete User account'
$ws.Range('B5').Value = 'Admin Controller'
$ws.Range('A6').Value = 'Promote User to Admin'
$ws.Range('B6').Value = 'Admin Controller'
$ws.Range('A8').Value = 'USER FEATURES'
$ws.Range('B8').Value = 'Area / Component'
$ws.Range('C8').Value = 'Status'
$ws.Range('D8').Value = 'Notes'
$ws.Range('A9').Value = 'Register'
$ws.Range('B9').Value = 'User Controller'
$ws.Range('C9').Value = 'Done'
$ws.Range('D9').Value = 'May need future updates to support faction avatar. '
$ws.Range('A10').Value = 'Login'
$ws.Range('B10').Value = 'User Controller'
$ws.Range('C10').Value = 'Done'
$ws.Range('A11').Value = 'Logout'
$ws.Range('B11').Value = 'User Controller'
$ws.Range('C11').Value = 'Done'
$ws.Range('A12').Value = 'Change Password'
$ws.Range('B12').Value = 'User Controller'
$ws.Range('C12').Value = 'Done'
$ws.Range('A13').Value = 'Change Avatar'
$ws.Range('B13').Value = 'User Controller'
$ws.Range('A14').Value = 'View User Profile'
$ws.Range('B14').Value = 'User Controller'
$ws.Range('A15').Value = 'View Users List'
$ws.Range('B15').Value = 'User Controller'
$ws.Range('A17').Value = 'MISC'
$ws.Range('B17').Value = 'Area / Component'
$ws.Range('C17').Value = 'Status'
$ws.Range('D17').Value = 'Notes'
$ws.Range('A18').Value = 'Allow OAUTH2 Logins from Facebook'
$ws.Range('B18').Value = 'Login/Auth'
$ws.Range('D18').Value = 'See StartupAuth.cs'
$ws.Range('A19').Value = 'Allow OAUTH2 Logins from Google'
$ws.Range('B19').Value = 'Login/Auth'
$ws.Range('D19').Value = 'See StartupAuth.cs'
$ws.Range('A20').Value = 'Modify/Update user model to include missing fields'
$ws.Range('B20').Value = 'User Model'
$ws.Range('C20').Value = 'Done'
$ws.Range('D20').Value = 'May need future updates to support faction avatar. '
$ws.Range('A21').Value = 'Add any missing/required fields from model to registration page'
$ws.Range('B21').Value = 'Registration View'
$ws.Range('C21').Value = 'Done'
$ws.Range('A22').Value = 'Add default avatars to database'
$ws.Range('A23').Value = 'Obscure/encrypt connection string ID/Password'
$ws.Range('B23').Value = 'Security'
$ws.Range('D23').Value = 'Connection string is in Web.Config'
$ws.Range('A24').Value = 'Fix code-first migrations on server'

# Section-header fills (copied from the matching header style on M-1 Tasks)
$styleSrc.Range('A1:D1').Copy()
$ws.Range('A17:D17').PasteSpecial(-4122)
$styleSrc.Range('A9:D9').Copy()
$ws.Range('A1:D1').PasteSpecial(-4122)
$styleSrc.Range('A16:D16').Copy()
$ws.Range('A8:D8').PasteSpecial(-4122)

# 'Done' status cells -> built-in Good (green) cell style
$ws.Range('C3').Style = 'Good'
$ws.Range('C9').Style = 'Good'
$ws.Range('C10').Style = 'Good'
$ws.Range('C11').Style = 'Good'
$ws.Range('C12').Style = 'Good'
$ws.Range('C20').Style = 'Good'
$ws.Range('C21').Style = 'Good'

# ===================== M1 - Game Data =====================
$ws = $wsM1

# Column widths
$ws.Columns.Item(1).ColumnWidth = 64.0
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 10.333333333333334
$ws.Columns.Item(4).ColumnWidth = 93.33333333333333

# Cell text
$ws.Range('A1').Value = 'ADMIN FEATURES'
$ws.Range('B1').Value = 'Area / Component'
$ws.Range('C1').Value = 'Status'
$ws.Range('D1').Value = 'Notes'
$ws.Range('A7').Value = 'USER FEATURES'
$ws.Range('B7').Value = 'Area / Component'
$ws.Range('C7').Value = 'Status'
$ws.Range('D7').Value = 'Notes'
$ws.Range('A13').Value = 'MISC'
$ws.Range('B13').Value = 'Area / Component'
$ws.Range('C13').Value = 'Status'
$ws.Range('D13').Value = 'Notes'

# Section-header fills (copied from the matching header style on M-1 Tasks)
$styleSrc.Range('A1:D1').Copy()
$ws.Range('A13:D13').PasteSpecial(-4122)
$styleSrc.Range('A9:D9').Copy()
$ws.Range('A1:D1').PasteSpecial(-4122)
$styleSrc.Range('A16:D16').Copy()
$ws.Range('A7:D7').PasteSpecial(-4122)

# ---- Selections on each sheet (left-over cursor position, matches target) ----
$wsM0.Range('A25').Select()
$wsM1.Range('A27').Select()

# M-0.5's cursor moved too (no longer the active tab)
$wb.Worksheets.Item('M-0.5').Range('C6').Select()

# M0 - Account Mgmt becomes the active/selected tab (activeTab=2, 0-based)
$wsM0.Activate()
$wsM0.Range('A25').Select()

Write-Output "Sheets now:"
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Output "$i : $($wb.Worksheets.Item($i).Name)"
}
